# Scheduled runner update: refresh market-board profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit* columns H:N) across all eight Sheets in the workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6 (Leve Item ID 4564)
$ws.Range("H6").Value = 4926508.5
$ws.Range("I6").Value = 6896872
$ws.Range("J6").Value = 600
$ws.Range("K6").Value = 20690616
$ws.Range("L6").Value = 1800
$ws.Range("M6").Value = -20690504
$ws.Range("N6").Value = -2024

# Row 8 (Leve Item ID 4565)
$ws.Range("H8").Value = 275.5102
$ws.Range("I8").Value = 181.28572
$ws.Range("J8").Value = 291.2143
$ws.Range("K8").Value = 543.85716
$ws.Range("L8").Value = 873.6428999999999
$ws.Range("M8").Value = -404.85716
$ws.Range("N8").Value = -1151.6429

# Row 29 (Leve Item ID 4575)
$ws.Range("H29").Value = 1750
$ws.Range("J29").Value = 2000
$ws.Range("L29").Value = 6000
$ws.Range("N29").Value = -6562

# Row 33 (Leve Item ID 5512)
$ws.Range("H33").Value = 17242364
$ws.Range("I33").Value = 21740072
$ws.Range("J33").Value = 1147.6666
$ws.Range("K33").Value = 21740072
$ws.Range("L33").Value = 1147.6666
$ws.Range("M33").Value = -21739843
$ws.Range("N33").Value = -1605.6666

# Row 38 (Leve Item ID 4599)
$ws.Range("H38").Value = 71.5
$ws.Range("I38").Value = 71.5
$ws.Range("K38").Value = 214.5
$ws.Range("M38").Value = 157.5

# Row 52 (Leve Item ID 4567)
$ws.Range("H52").Value = 299
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

# Row 62 (Leve Item ID 27781)
$ws.Range("H62").Value = 88881.2
$ws.Range("J62").Value = 88881.2
$ws.Range("L62").Value = 88881.2
$ws.Range("N62").Value = -90129.2

# Row 65 (Leve Item ID 27781)
$ws.Range("H65").Value = 88881.2
$ws.Range("J65").Value = 88881.2
$ws.Range("L65").Value = 444406
$ws.Range("N65").Value = -450646

# Row 96 (Leve Item ID 19894)
$ws.Range("H96").Value = 700.3333
$ws.Range("I96").Value = 532.25
$ws.Range("K96").Value = 1596.75
$ws.Range("M96").Value = -223.75

# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 1445
$ws.Range("I132").Value = 1341
$ws.Range("K132").Value = 4023
$ws.Range("M132").Value = -1493

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 2685.1372
$ws.Range("I137").Value = 2468.054
$ws.Range("K137").Value = 7404.162
$ws.Range("M137").Value = -4854.162

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 5697.9287
$ws.Range("J138").Value = 10469.096
$ws.Range("L138").Value = 31407.288
$ws.Range("N138").Value = -41687.288

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 3995.52
$ws.Range("I32").Value = 3822.032
$ws.Range("K32").Value = 3822.032
$ws.Range("M32").Value = -3535.032

# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 5311.025
$ws.Range("I61").Value = 1478.72
$ws.Range("J61").Value = 11698.2
$ws.Range("K61").Value = 1478.72
$ws.Range("L61").Value = 11698.2
$ws.Range("M61").Value = -1266.72
$ws.Range("N61").Value = -12122.2

# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 2163.4688
$ws.Range("I74").Value = 1161.7727
$ws.Range("J74").Value = 4367.2
$ws.Range("K74").Value = 1161.7727
$ws.Range("L74").Value = 4367.2
$ws.Range("M74").Value = -287.7727
$ws.Range("N74").Value = -6115.2

# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 2163.4688
$ws.Range("I77").Value = 1161.7727
$ws.Range("J77").Value = 4367.2
$ws.Range("K77").Value = 5808.863499999999
$ws.Range("L77").Value = 21836
$ws.Range("M77").Value = -1440.863499999999
$ws.Range("N77").Value = -30572

# Row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 17545112
$ws.Range("I110").Value = 1168
$ws.Range("K110").Value = 1168
$ws.Range("M110").Value = 877

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 5311.025
$ws.Range("I136").Value = 1478.72
$ws.Range("J136").Value = 11698.2
$ws.Range("K136").Value = 4436.16
$ws.Range("L136").Value = 35094.60000000001
$ws.Range("M136").Value = -1886.16
$ws.Range("N136").Value = -40194.60000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (Leve Item ID 14149)
$ws.Range("H20").Value = 5556969.5
$ws.Range("I20").Value = 8334674
$ws.Range("K20").Value = 8334674
$ws.Range("M20").Value = -8334427

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 4362.557
$ws.Range("I134").Value = 1704.4667
$ws.Range("J134").Value = 9147.120000000001
$ws.Range("K134").Value = 5113.4001
$ws.Range("L134").Value = 27441.36
$ws.Range("M134").Value = -2578.4001
$ws.Range("N134").Value = -32511.36

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 6053.0747
$ws.Range("I31").Value = 2913
$ws.Range("J31").Value = 10705.037
$ws.Range("K31").Value = 2913
$ws.Range("L31").Value = 10705.037
$ws.Range("M31").Value = -2618
$ws.Range("N31").Value = -11295.037

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 6053.0747
$ws.Range("I34").Value = 2913
$ws.Range("J34").Value = 10705.037
$ws.Range("K34").Value = 2913
$ws.Range("L34").Value = 10705.037
$ws.Range("M34").Value = -2711
$ws.Range("N34").Value = -11109.037

# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 10644031
$ws.Range("J58").Value = 9768.174000000001
$ws.Range("L58").Value = 9768.174000000001
$ws.Range("N58").Value = -10174.174

# Row 104 (Leve Item ID 19749)
$ws.Range("H104").Value = 39285
$ws.Range("J104").Value = 39285
$ws.Range("L104").Value = 39285
$ws.Range("N104").Value = -44527

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 4880.705
$ws.Range("I132").Value = 2575.1396
$ws.Range("K132").Value = 7725.418799999999
$ws.Range("M132").Value = -5195.418799999999

# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 3795.4546
$ws.Range("I134").Value = 1664.2982
$ws.Range("K134").Value = 4992.8946
$ws.Range("M134").Value = -2457.8946

# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 10644031
$ws.Range("J136").Value = 9768.174000000001
$ws.Range("L136").Value = 29304.522
$ws.Range("N136").Value = -34404.522

$ws = $wb.Worksheets.Item("CUL")
# Row 129 (Leve Item ID 36054)
$ws.Range("H129").Value = 850.7273
$ws.Range("I129").Value = 536.625
$ws.Range("J129").Value = 1688.3334
$ws.Range("K129").Value = 1609.875
$ws.Range("L129").Value = 5065.0002
$ws.Range("M129").Value = 3390.125
$ws.Range("N129").Value = -15065.0002

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (Leve Item ID 5062)
$ws.Range("H2").Value = 198.21428
$ws.Range("I2").Value = 63.666668
$ws.Range("J2").Value = 299.125
$ws.Range("K2").Value = 63.666668
$ws.Range("L2").Value = 299.125
$ws.Range("M2").Value = 49.333332
$ws.Range("N2").Value = -525.125

# Row 107 (Leve Item ID 27802)
$ws.Range("H107").Value = 889231.9
$ws.Range("I107").Value = 1333602.1
$ws.Range("J107").Value = 491.33334
$ws.Range("K107").Value = 1333602.1
$ws.Range("L107").Value = 491.33334
$ws.Range("M107").Value = -1331682.1
$ws.Range("N107").Value = -4331.33334

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 15875787
$ws.Range("I46").Value = 1475
$ws.Range("K46").Value = 1475
$ws.Range("M46").Value = -1287

# Row 100 (Leve Item ID 19995)
$ws.Range("H100").Value = 3831.8635
$ws.Range("I100").Value = 2894.1667
$ws.Range("J100").Value = 4957.1
$ws.Range("K100").Value = 2894.1667
$ws.Range("L100").Value = 4957.1
$ws.Range("M100").Value = -2353.1667
$ws.Range("N100").Value = -6039.1

# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 5461.75
$ws.Range("I122").Value = 4205.857
$ws.Range("J122").Value = 7220
$ws.Range("K122").Value = 12617.571
$ws.Range("L122").Value = 21660
$ws.Range("M122").Value = -10167.571
$ws.Range("N122").Value = -26560

# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 9265551
$ws.Range("I132").Value = 15628292
$ws.Range("J132").Value = 10656.546
$ws.Range("K132").Value = 46884876
$ws.Range("L132").Value = 31969.638
$ws.Range("M132").Value = -46882346
$ws.Range("N132").Value = -37029.638

$ws = $wb.Worksheets.Item("WVR")
# Row 20 (Leve Item ID 3023)
$ws.Range("H20").Value = 9000
$ws.Range("J20").Value = 9000
$ws.Range("L20").Value = 9000
$ws.Range("N20").Value = -9480

